$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.383.73"
$ws.Range("E2").Value = "  -1.31%  "
$ws.Range("D3").Value = "2.341.70"
$ws.Range("E3").Value = "  +3.04%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("B5").Value = "XRP"
$ws.Range("C5").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D5").Value = "'0.651"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.14%  "
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "'232.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("D7").Value = "'65.96"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.54%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "'0.453"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.88%  "
$ws.Range("D10").Value = "'0.0955"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.62%  "
$ws.Range("D11").Value = "'56.95"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.02%  "
$ws.Range("D12").Value = "'26.81"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.84%  "
$ws.Range("D13").Value = "2.688.78"
$ws.Range("E13").Value = "  +2.99%  "
$ws.Range("D14").Value = "'0.104"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("D15").Value = "'15.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.68%  "
$ws.Range("D16").Value = "'6.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.79%  "
$ws.Range("D17").Value = "'0.847"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.08%  "
$ws.Range("D18").Value = "2.342.83"
$ws.Range("E18").Value = "  +3.03%  "
$ws.Range("D19").Value = "43.323.76"
$ws.Range("E19").Value = "  -1.24%  "
$ws.Range("D20").Value = "0.0₃0979"
$ws.Range("E20").Value = "  -2.47%  "
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("D22").Value = "'6.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.41%  "
$ws.Range("D23").Value = "'248.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.55%  "
$ws.Range("D24").Value = "'3.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +15.33%  "
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("E26").Value = "  -0.96%  "
$ws.Range("E27").Value = "  -1.43%  "
$ws.Range("D28").Value = "'9.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.00%  "
$ws.Range("D29").Value = "'175.12"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.91%  "
$ws.Range("E30").Value = "  +6.06%  "
$ws.Range("D31").Value = "'1.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.60%  "
$ws.Range("E32").Value = "  -7.76%  "
$ws.Range("D33").Value = "'0.126"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.60%  "
$ws.Range("E34").Value = "  +4.06%  "
$ws.Range("E35").Value = "  -2.47%  "
$ws.Range("D36").Value = "'4.96"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.89%  "
$ws.Range("D37").Value = "'2.52"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.19%  "
$ws.Range("D38").Value = "'6.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("E39").Value = "  -5.34%  "
$ws.Range("E40").Value = "  -3.14%  "
$ws.Range("E41").Value = "  +8.93%  "
$ws.Range("E42").Value = "  -0.18%  "
$ws.Range("D43").Value = "'18.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.22%  "
$ws.Range("E44").Value = "  +9.10%  "
$ws.Range("D45").Value = "'99.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.75%  "
$ws.Range("E46").Value = "  -0.81%  "
$ws.Range("D47").Value = "'4.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.07%  "
$ws.Range("E48").Value = "  -4.34%  "
$ws.Range("D49").Value = "1.438.93"
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("B50").Value = "Celestia"
$ws.Range("C50").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D50").Value = "'9.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.05%  "
$ws.Range("D51").Value = "'0.000203"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -10.70%  "
